$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 36, pushing existing rows 36-38 down to 37-39.
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new weekly data point.
$ws.Cells.Item(36, 1).Value = 7
$ws.Cells.Item(36, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(36, 3).Value = "Ñuble"
$ws.Cells.Item(36, 4).Value = 44776
$ws.Cells.Item(36, 5).Value = 16
$ws.Cells.Item(36, 6).Value = 100112001
$ws.Cells.Item(36, 7).Value = "Berenjena"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 60
$ws.Cells.Item(36, 11).Value = 11000
$ws.Cells.Item(36, 12).Value = 12000
$ws.Cells.Item(36, 13).Value = 11500
$ws.Cells.Item(36, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(36, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(36, 16).Value = 192
$ws.Cells.Item(36, 17).Value = 60
$ws.Cells.Item(36, 18).Value = "Hortaliza"

# Match the date formatting style used by column D elsewhere (row 35 as reference).
$ws.Cells.Item(36, 4).NumberFormat = $ws.Cells.Item(35, 4).NumberFormat
